$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 147, pushing existing rows 147:209 down to 148:210
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new data record
$ws.Range("A147").Value = 3
$ws.Range("B147").Value = "Femacal de La Calera"
$ws.Range("C147").Value = "Coquimbo"
$ws.Range("D147").Value = 44489
$ws.Range("E147").Value = 5
$ws.Range("F147").Value = 100114013
$ws.Range("G147").Value = "Zanahoria"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 180
$ws.Range("K147").Value = 8000
$ws.Range("L147").Value = 8000
$ws.Range("M147").Value = 8000
$ws.Range("N147").Value = "`$/saco 20 kilos"
$ws.Range("O147").Value = "Provincia de Quillota"
$ws.Range("P147").Value = 400
$ws.Range("Q147").Value = 20
$ws.Range("R147").Value = "Hortaliza"
